$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-8, columns B (TB), C (d2S), D (K), E (IP), F (Win), G (sum)
$data = @{
    2 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 0, 10.35301142835362)
    3 = @(0.06328177979961902, 0.3375848360084654, 16.98373111632243, 6.48142807727062, 1, 23.86602580940113)
    4 = @(1.505614041169197, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 0, 3.055818435266709)
    5 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 3.034748368925986)
    6 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 6.741336633845642)
    7 = @(3.182878228561681, 9.226618575922256, 0.7127328510149897, 6.48142807727062, 0, 19.60365773276954)
    8 = @(3.182878228561681, 86.29678392075563, 16.98373111632243, 6.48142807727062, 0, 112.9448213429104)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # column B is 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
